$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Prepare formatting for the three "new" row positions (21, 22, 23)
#    by copying formats from the rows that currently hold that look,
#    BEFORE any values are overwritten.
# ------------------------------------------------------------------

# Row 22 (old footer) -> Row 23 (new footer position)
$ws.Range("A22:N22").Copy()
$ws.Range("A23:N23").PasteSpecial(-4122)

# Row 21 (old totals) -> Row 22 (new totals position)
$ws.Range("K21:N21").Copy()
$ws.Range("K22:N22").PasteSpecial(-4122)

# Row 20 (last item row) -> Row 21 (new item row)
$ws.Range("A20:N20").Copy()
$ws.Range("A21:N21").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2) Shift the item rows down by one: row 20 -> 21, 19 -> 20, ... 10 -> 11
#    Only B (name), H (balance text), L (price), N (ratio) vary per row;
#    A just keeps the running 1..18 count already in place.
# ------------------------------------------------------------------
for ($r = 20; $r -ge 10; $r--) {
    $dst = $r + 1
    $ws.Range("B$dst").Value2 = $ws.Range("B$r").Value2
    $ws.Range("H$dst").Value2 = $ws.Range("H$r").Value2
    $ws.Range("L$dst").Value2 = $ws.Range("L$r").Value2
    $ws.Range("N$dst").Value2 = $ws.Range("N$r").Value2
}

# ------------------------------------------------------------------
# 3) Write the new item into row 10 (alphabetically between COLOVATIL
#    and GAVISCON): DOSIN 2MG 20 TAB.
# ------------------------------------------------------------------
$ws.Range("B10").Value2 = "DOSIN 2MG 20 TAB."
$ws.Range("H10").Value2 = "1:0"
$ws.Range("L10").Value2 = 32
$ws.Range("N10").Value2 = 1

# ------------------------------------------------------------------
# 4) Totals row now lives at row 22; update its value.
# ------------------------------------------------------------------
$ws.Range("K22").Value2 = 1238.35

# ------------------------------------------------------------------
# 5) Footer row now lives at row 23; restore its text.
# ------------------------------------------------------------------
$ws.Range("A23").Value2 = "Monday, 5 January, 2026 10:41 AM"
$ws.Range("F23").Value2 = "1/1"
$ws.Range("I23").Value2 = "developed by : Abdelaziz Talaat"

# ------------------------------------------------------------------
# 6) Row heights: new item row 21 matches the alternating pattern,
#    totals row 22 grows slightly, footer row 23 unchanged.
# ------------------------------------------------------------------
$ws.Rows("21:21").RowHeight = 24.75
$ws.Rows("22:22").RowHeight = 26.25
$ws.Rows("23:23").RowHeight = 16.5

# ------------------------------------------------------------------
# 7) Merged cells for the new row layout.
# ------------------------------------------------------------------
$ws.Range("B21:G21").Merge()
$ws.Range("H21:K21").Merge()
$ws.Range("L21:M21").Merge()
$ws.Range("K22:N22").Merge()
$ws.Range("A23:E23").Merge()
$ws.Range("F23:G23").Merge()
$ws.Range("I23:N23").Merge()
